$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.165.20"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").Value = "3.746.25"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'614.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").Value = "'178.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").Value = "3.744.63"
$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -1.80%  "

$ws.Range("D10").Value = "'0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.58%  "

$ws.Range("E11").Value = "  +2.59%  "

$ws.Range("D12").Value = "'0.484"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.34%  "

$ws.Range("D13").Value = "'40.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.42%  "

$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").Value = "4.366.04"
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").Value = "3.745.23"
$ws.Range("E16").Value = "  -1.09%  "

$ws.Range("D17").Value = "69.247.27"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D19").Value = "'7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.00%  "

$ws.Range("D20").Value = "'16.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.86%  "

$ws.Range("D21").Value = "'498.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.24%  "

$ws.Range("E22").Value = "  -2.61%  "

$ws.Range("D23").Value = "'0.721"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D24").Value = "'2.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("D25").Value = "'85.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.79%  "

$ws.Range("D26").Value = "'12.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.59%  "

$ws.Range("D27").Value = "'10.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.53%  "

$ws.Range("E28").Value = "  -3.81%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").Value = "'2.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "

$ws.Range("E31").Value = "  +3.31%  "

$ws.Range("D32").Value = "'8.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.01%  "

$ws.Range("D33").Value = "'30.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.36%  "

$ws.Range("E34").Value = "  -1.87%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("E36").Value = "  -0.39%  "

$ws.Range("D37").Value = "'6.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.82%  "

$ws.Range("D38").Value = "'0.347"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.76%  "

$ws.Range("D40").Value = "'450.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.85%  "

$ws.Range("E41").Value = "  +8.94%  "

$ws.Range("E42").Value = "  -5.39%  "

$ws.Range("D43").Value = "'49.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.17%  "

$ws.Range("E44").Value = "  +0.68%  "

$ws.Range("D45").Value = "'8.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.68%  "

$ws.Range("D46").Value = "2.946.76"
$ws.Range("E46").Value = "  -4.00%  "

$ws.Range("D47").Value = "'0.0359"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.94%  "

$ws.Range("D48").Value = "'27.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.34%  "

$ws.Range("D50").Value = "'137.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.47%  "

$ws.Range("E51").Value = "  -1.45%  "
